$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.858.96"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "1.894.43"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7777"
$ws.Range("E5").Value = "  -4.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.80"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3148"
$ws.Range("E8").Value = "  -3.28%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07394"
$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.51"
$ws.Range("E10").Value = "  -5.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08106"
$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7730"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.520"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "1.840.93"
$ws.Range("E14").Value = "  -3.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.47"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.323"
$ws.Range("E16").Value = "  +4.20%  "

$ws.Range("D17").Value = "29.802.81"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.00"
$ws.Range("E18").Value = "  -1.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.25"
$ws.Range("E19").Value = "  -1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007837"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.153"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "2.128.94"
$ws.Range("E23").Value = "  -2.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1590"
$ws.Range("E25").Value = "  -5.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.482"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.72"
$ws.Range("E27").Value = "  -3.10%  "

$ws.Range("E28").Value = "  -1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.050"
$ws.Range("E29").Value = "  -5.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  +4.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.553"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.488"
$ws.Range("E32").Value = "  +3.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.103"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05542"
$ws.Range("E34").Value = "  -5.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("E36").Value = "  +1.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.639"
$ws.Range("E38").Value = "  -3.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01927"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").Value = "1.159.50"
$ws.Range("E41").Value = "  +13.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.92"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4441"
$ws.Range("E43").Value = "  -2.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.945"
$ws.Range("E44").Value = "  -0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8501"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.897"
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.35"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.088"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.501"
$ws.Range("E51").Value = "  -1.54%  "
